# update redirects after successful ad change. update for ad edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Deactivate/Publish Again Ad" (row 25) and "Delete Ad" (row 28) are now
# implemented, so they score full marks - Total Score (C51) recalculates
# automatically via its SUM formula.
$ws.Range("C25").Value = 5
$ws.Range("C28").Value = 5

# Reflect the scrolled/selected state left after making the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("H18").Select()
